$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 271, shifting existing rows 271-309 down to 272-310
$ws.Rows.Item(271).EntireRow.Insert()

# Populate the newly inserted row 271 with the new weekly price entry
$ws.Range("A271").Value = 4
$ws.Range("B271").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C271").Value = "Los Lagos"
$ws.Range("D271").Value = 44776
$ws.Range("E271").Value = 10
$ws.Range("F271").Value = 100112037
$ws.Range("G271").Value = "Cebollín"
$ws.Range("H271").Value = "Sin especificar"
$ws.Range("I271").Value = "Primera"
$ws.Range("J271").Value = 35
$ws.Range("K271").Value = 10000
$ws.Range("L271").Value = 10000
$ws.Range("M271").Value = 10000
$ws.Range("N271").Value = "$/paquete 36 unidades"
$ws.Range("O271").Value = "Región Metropolitana"
$ws.Range("P271").Value = 278
$ws.Range("Q271").Value = 36
$ws.Range("R271").Value = "Hortaliza"
